$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.249.38"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.268.88"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.495"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.33%  "
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "2.619.50"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "2.269.14"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.785"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").Value = "42.137.96"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0686"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").Value = "1.973.09"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0278"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.20%  "
$ws.Range("E47").Value = "  -5.57%  "
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("D49").Value = "2.492.44"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.27%  "
